$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the original values for columns D,L,M,N,O,P,Q,R,S,T for rows 2-19
# (columns A-C, E-K are constant across rows and are left untouched).
$orig = @{}
$orig[2] = @{
  "D" = $ws.Range("D2").Value2
  "L" = $ws.Range("L2").Value2
  "M" = $ws.Range("M2").Value2
  "N" = $ws.Range("N2").Value2
  "O" = $ws.Range("O2").Value2
  "P" = $ws.Range("P2").Value2
  "Q" = $ws.Range("Q2").Value2
  "R" = $ws.Range("R2").Value2
  "S" = $ws.Range("S2").Value2
  "T" = $ws.Range("T2").Value2
}
$orig[3] = @{
  "D" = $ws.Range("D3").Value2
  "L" = $ws.Range("L3").Value2
  "M" = $ws.Range("M3").Value2
  "N" = $ws.Range("N3").Value2
  "O" = $ws.Range("O3").Value2
  "P" = $ws.Range("P3").Value2
  "Q" = $ws.Range("Q3").Value2
  "R" = $ws.Range("R3").Value2
  "S" = $ws.Range("S3").Value2
  "T" = $ws.Range("T3").Value2
}
$orig[4] = @{
  "D" = $ws.Range("D4").Value2
  "L" = $ws.Range("L4").Value2
  "M" = $ws.Range("M4").Value2
  "N" = $ws.Range("N4").Value2
  "O" = $ws.Range("O4").Value2
  "P" = $ws.Range("P4").Value2
  "Q" = $ws.Range("Q4").Value2
  "R" = $ws.Range("R4").Value2
  "S" = $ws.Range("S4").Value2
  "T" = $ws.Range("T4").Value2
}
$orig[5] = @{
  "D" = $ws.Range("D5").Value2
  "L" = $ws.Range("L5").Value2
  "M" = $ws.Range("M5").Value2
  "N" = $ws.Range("N5").Value2
  "O" = $ws.Range("O5").Value2
  "P" = $ws.Range("P5").Value2
  "Q" = $ws.Range("Q5").Value2
  "R" = $ws.Range("R5").Value2
  "S" = $ws.Range("S5").Value2
  "T" = $ws.Range("T5").Value2
}
$orig[6] = @{
  "D" = $ws.Range("D6").Value2
  "L" = $ws.Range("L6").Value2
  "M" = $ws.Range("M6").Value2
  "N" = $ws.Range("N6").Value2
  "O" = $ws.Range("O6").Value2
  "P" = $ws.Range("P6").Value2
  "Q" = $ws.Range("Q6").Value2
  "R" = $ws.Range("R6").Value2
  "S" = $ws.Range("S6").Value2
  "T" = $ws.Range("T6").Value2
}
$orig[7] = @{
  "D" = $ws.Range("D7").Value2
  "L" = $ws.Range("L7").Value2
  "M" = $ws.Range("M7").Value2
  "N" = $ws.Range("N7").Value2
  "O" = $ws.Range("O7").Value2
  "P" = $ws.Range("P7").Value2
  "Q" = $ws.Range("Q7").Value2
  "R" = $ws.Range("R7").Value2
  "S" = $ws.Range("S7").Value2
  "T" = $ws.Range("T7").Value2
}
$orig[8] = @{
  "D" = $ws.Range("D8").Value2
  "L" = $ws.Range("L8").Value2
  "M" = $ws.Range("M8").Value2
  "N" = $ws.Range("N8").Value2
  "O" = $ws.Range("O8").Value2
  "P" = $ws.Range("P8").Value2
  "Q" = $ws.Range("Q8").Value2
  "R" = $ws.Range("R8").Value2
  "S" = $ws.Range("S8").Value2
  "T" = $ws.Range("T8").Value2
}
$orig[9] = @{
  "D" = $ws.Range("D9").Value2
  "L" = $ws.Range("L9").Value2
  "M" = $ws.Range("M9").Value2
  "N" = $ws.Range("N9").Value2
  "O" = $ws.Range("O9").Value2
  "P" = $ws.Range("P9").Value2
  "Q" = $ws.Range("Q9").Value2
  "R" = $ws.Range("R9").Value2
  "S" = $ws.Range("S9").Value2
  "T" = $ws.Range("T9").Value2
}
$orig[10] = @{
  "D" = $ws.Range("D10").Value2
  "L" = $ws.Range("L10").Value2
  "M" = $ws.Range("M10").Value2
  "N" = $ws.Range("N10").Value2
  "O" = $ws.Range("O10").Value2
  "P" = $ws.Range("P10").Value2
  "Q" = $ws.Range("Q10").Value2
  "R" = $ws.Range("R10").Value2
  "S" = $ws.Range("S10").Value2
  "T" = $ws.Range("T10").Value2
}
$orig[11] = @{
  "D" = $ws.Range("D11").Value2
  "L" = $ws.Range("L11").Value2
  "M" = $ws.Range("M11").Value2
  "N" = $ws.Range("N11").Value2
  "O" = $ws.Range("O11").Value2
  "P" = $ws.Range("P11").Value2
  "Q" = $ws.Range("Q11").Value2
  "R" = $ws.Range("R11").Value2
  "S" = $ws.Range("S11").Value2
  "T" = $ws.Range("T11").Value2
}
$orig[12] = @{
  "D" = $ws.Range("D12").Value2
  "L" = $ws.Range("L12").Value2
  "M" = $ws.Range("M12").Value2
  "N" = $ws.Range("N12").Value2
  "O" = $ws.Range("O12").Value2
  "P" = $ws.Range("P12").Value2
  "Q" = $ws.Range("Q12").Value2
  "R" = $ws.Range("R12").Value2
  "S" = $ws.Range("S12").Value2
  "T" = $ws.Range("T12").Value2
}
$orig[13] = @{
  "D" = $ws.Range("D13").Value2
  "L" = $ws.Range("L13").Value2
  "M" = $ws.Range("M13").Value2
  "N" = $ws.Range("N13").Value2
  "O" = $ws.Range("O13").Value2
  "P" = $ws.Range("P13").Value2
  "Q" = $ws.Range("Q13").Value2
  "R" = $ws.Range("R13").Value2
  "S" = $ws.Range("S13").Value2
  "T" = $ws.Range("T13").Value2
}
$orig[14] = @{
  "D" = $ws.Range("D14").Value2
  "L" = $ws.Range("L14").Value2
  "M" = $ws.Range("M14").Value2
  "N" = $ws.Range("N14").Value2
  "O" = $ws.Range("O14").Value2
  "P" = $ws.Range("P14").Value2
  "Q" = $ws.Range("Q14").Value2
  "R" = $ws.Range("R14").Value2
  "S" = $ws.Range("S14").Value2
  "T" = $ws.Range("T14").Value2
}
$orig[15] = @{
  "D" = $ws.Range("D15").Value2
  "L" = $ws.Range("L15").Value2
  "M" = $ws.Range("M15").Value2
  "N" = $ws.Range("N15").Value2
  "O" = $ws.Range("O15").Value2
  "P" = $ws.Range("P15").Value2
  "Q" = $ws.Range("Q15").Value2
  "R" = $ws.Range("R15").Value2
  "S" = $ws.Range("S15").Value2
  "T" = $ws.Range("T15").Value2
}
$orig[16] = @{
  "D" = $ws.Range("D16").Value2
  "L" = $ws.Range("L16").Value2
  "M" = $ws.Range("M16").Value2
  "N" = $ws.Range("N16").Value2
  "O" = $ws.Range("O16").Value2
  "P" = $ws.Range("P16").Value2
  "Q" = $ws.Range("Q16").Value2
  "R" = $ws.Range("R16").Value2
  "S" = $ws.Range("S16").Value2
  "T" = $ws.Range("T16").Value2
}
$orig[17] = @{
  "D" = $ws.Range("D17").Value2
  "L" = $ws.Range("L17").Value2
  "M" = $ws.Range("M17").Value2
  "N" = $ws.Range("N17").Value2
  "O" = $ws.Range("O17").Value2
  "P" = $ws.Range("P17").Value2
  "Q" = $ws.Range("Q17").Value2
  "R" = $ws.Range("R17").Value2
  "S" = $ws.Range("S17").Value2
  "T" = $ws.Range("T17").Value2
}
$orig[18] = @{
  "D" = $ws.Range("D18").Value2
  "L" = $ws.Range("L18").Value2
  "M" = $ws.Range("M18").Value2
  "N" = $ws.Range("N18").Value2
  "O" = $ws.Range("O18").Value2
  "P" = $ws.Range("P18").Value2
  "Q" = $ws.Range("Q18").Value2
  "R" = $ws.Range("R18").Value2
  "S" = $ws.Range("S18").Value2
  "T" = $ws.Range("T18").Value2
}
$orig[19] = @{
  "D" = $ws.Range("D19").Value2
  "L" = $ws.Range("L19").Value2
  "M" = $ws.Range("M19").Value2
  "N" = $ws.Range("N19").Value2
  "O" = $ws.Range("O19").Value2
  "P" = $ws.Range("P19").Value2
  "Q" = $ws.Range("Q19").Value2
  "R" = $ws.Range("R19").Value2
  "S" = $ws.Range("S19").Value2
  "T" = $ws.Range("T19").Value2
}

# Row permutation observed in the diff: row $dest ends up holding the data
# that row $src held before the edit.
$mapping = @{
  2 = 14
  3 = 15
  4 = 18
  5 = 4
  6 = 5
  7 = 2
  8 = 6
  9 = 7
  10 = 13
  11 = 17
  12 = 19
  13 = 11
  14 = 9
  15 = 16
  16 = 10
  17 = 8
  18 = 3
  19 = 12
}

foreach ($dest in $mapping.Keys) {
  $src = $mapping[$dest]
  $row = $orig[$src]
  $ws.Range("D$dest").Value2 = $row["D"]
  $ws.Range("L$dest").Value2 = $row["L"]
  $ws.Range("M$dest").Value2 = $row["M"]
  $ws.Range("N$dest").Value2 = $row["N"]
  $ws.Range("O$dest").Value2 = $row["O"]
  $ws.Range("P$dest").Value2 = $row["P"]
  $ws.Range("Q$dest").Value2 = $row["Q"]
  $ws.Range("R$dest").Value2 = $row["R"]
  $ws.Range("S$dest").Value2 = $row["S"]
  $ws.Range("T$dest").Value2 = $row["T"]
}
